$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 98
$ws.Range("J2").Value = 434
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 106
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 65
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 36
$ws.Range("T2").Value = 74
$ws.Range("U2").Value = 1
$ws.Range("V2").Value = 667
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 666
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 8
$ws.Range("AA2").Value = 1
